$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the old "mation"/"pompes)"/"Hiver"/"Eté"/"Année" junk header row (row 2),
# shifting the data rows (old 3..7) up to become rows 2..6.
$ws.Rows.Item(2).Delete()

# Re-label the header row. Columns A-E are brand new index/name/date columns
# (plain, unstyled); F keeps the existing "(m3/s)" header while G-K become the
# new MW/GWh headers - all get the Arial 9pt header font, with General number
# format (no explicit numFmt application).
$tmpStyle = $wb.Styles.Add("__tmp_header_style")
$tmpStyle.Font.Name = "Arial"
$tmpStyle.Font.Size = 9

$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("E1").ClearFormats()

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"
$ws.Range("F1:K1").Style = "__tmp_header_style"

# Drop the scratch named style again - the cellXfs entry it produced survives
# on the cells (xfId falls back to the default "Normal" cellStyleXf) but the
# extra named cellStyle bookkeeping is removed.
$wb.Styles.Item("__tmp_header_style").Delete()

# New selection, matching the post-edit workbook (active cell A2, row 2 selected).
$ws.Range("A2:K2").Select()

$wb.Save()
